$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 936.4286
$ws.Range("J17").Value = 1042.75
$ws.Range("L17").Value = 3128.25
$ws.Range("N17").Value = -3464.25

$ws.Range("H98").Value = 1329.8
$ws.Range("I98").Value = 1329.8
$ws.Range("K98").Value = 1329.8
$ws.Range("M98").Value = 168.2

$ws.Range("H106").Value = 4309.0713
$ws.Range("I106").Value = 4309.0713
$ws.Range("K106").Value = 4309.0713
$ws.Range("M106").Value = -3678.0713

$ws.Range("H122").Value = 1329.8
$ws.Range("I122").Value = 1329.8
$ws.Range("K122").Value = 3989.4
$ws.Range("M122").Value = -1539.4

$ws.Range("H129").Value = 1974.95
$ws.Range("I129").Value = 1631.375
$ws.Range("K129").Value = 4894.125
$ws.Range("M129").Value = 105.875

$ws.Range("H131").Value = 20001580
$ws.Range("I131").Value = 25000700
$ws.Range("K131").Value = 75002100
$ws.Range("M131").Value = -74997060

$ws.Range("H132").Value = 1667.5714
$ws.Range("I132").Value = 1370.262
$ws.Range("J132").Value = 3451.4285
$ws.Range("K132").Value = 4110.786
$ws.Range("L132").Value = 10354.2855
$ws.Range("M132").Value = -1580.786
$ws.Range("N132").Value = -15414.2855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5454.0625
$ws.Range("I32").Value = 2385.453
$ws.Range("K32").Value = 2385.453
$ws.Range("M32").Value = -2098.453

$ws.Range("H61").Value = 42141.32
$ws.Range("I61").Value = 1986.7
$ws.Range("K61").Value = 1986.7
$ws.Range("M61").Value = -1774.7

$ws.Range("H102").Value = 47802.375
$ws.Range("I102").Value = 51173.5
$ws.Range("J102").Value = 30946.75
$ws.Range("K102").Value = 51173.5
$ws.Range("L102").Value = 30946.75
$ws.Range("M102").Value = -49551.5
$ws.Range("N102").Value = -34190.75

$ws.Range("H136").Value = 42141.32
$ws.Range("I136").Value = 1986.7
$ws.Range("K136").Value = 5960.1
$ws.Range("M136").Value = -3410.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1219.6842
$ws.Range("I20").Value = 1109
$ws.Range("K20").Value = 1109
$ws.Range("M20").Value = -862

$ws.Range("H105").Value = 42867.32
$ws.Range("I105").Value = 54875.367
$ws.Range("K105").Value = 54875.367
$ws.Range("M105").Value = -53128.367

$ws.Range("H134").Value = 5059.081
$ws.Range("I134").Value = 2258.3333
$ws.Range("K134").Value = 6774.999899999999
$ws.Range("M134").Value = -4239.999899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()

$ws.Range("H22").Value = 1318.1818
$ws.Range("J22").Value = 1275
$ws.Range("L22").Value = 1275
$ws.Range("N22").Value = -1975

$ws.Range("H58").Value = 1544.8572
$ws.Range("I58").Value = 1468.7646
$ws.Range("K58").Value = 1468.7646
$ws.Range("M58").Value = -1265.7646

$ws.Range("H122").Value = 3063.45
$ws.Range("I122").Value = 1969.6666
$ws.Range("J122").Value = 3958.3635
$ws.Range("K122").Value = 5908.9998
$ws.Range("L122").Value = 11875.0905
$ws.Range("M122").Value = -3458.9998
$ws.Range("N122").Value = -16775.0905

$ws.Range("H132").Value = 1301022.9
$ws.Range("I132").Value = 1569530.2
$ws.Range("K132").Value = 4708590.6
$ws.Range("M132").Value = -4706060.6

$ws.Range("H134").Value = 2132521.2
$ws.Range("J134").Value = 112535.22
$ws.Range("L134").Value = 337605.66
$ws.Range("N134").Value = -342675.66

$ws.Range("H136").Value = 1544.8572
$ws.Range("I136").Value = 1468.7646
$ws.Range("K136").Value = 4406.293799999999
$ws.Range("M136").Value = -1856.293799999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8121.778
$ws.Range("I70").Value = 8849.333000000001
$ws.Range("J70").Value = 6666.6665
$ws.Range("K70").Value = 8849.333000000001
$ws.Range("L70").Value = 6666.6665
$ws.Range("M70").Value = -8579.333000000001
$ws.Range("N70").Value = -7206.6665

$ws.Range("H73").Value = 8121.778
$ws.Range("I73").Value = 8849.333000000001
$ws.Range("J73").Value = 6666.6665
$ws.Range("K73").Value = 8849.333000000001
$ws.Range("L73").Value = 6666.6665
$ws.Range("M73").Value = -7913.333000000001
$ws.Range("N73").Value = -8538.666499999999

$ws.Range("H102").Value = 2337.4
$ws.Range("I102").Value = 2239.4285
$ws.Range("K102").Value = 2239.4285
$ws.Range("M102").Value = -617.4285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 68666.664
$ws.Range("J6").Value = 68666.664
$ws.Range("L6").Value = 68666.664
$ws.Range("N6").Value = -68890.664

$ws.Range("H40").Value = 10373490
$ws.Range("I40").Value = 3432.889
$ws.Range("K40").Value = 3432.889
$ws.Range("M40").Value = -3296.889

$ws.Range("H82").Value = 2199.75
$ws.Range("I82").Value = 1933
$ws.Range("K82").Value = 1933
$ws.Range("M82").Value = -1572

$ws.Range("H85").Value = 2199.75
$ws.Range("I85").Value = 1933
$ws.Range("K85").Value = 1933
$ws.Range("M85").Value = -685

$ws.Range("H136").Value = 1739.919
$ws.Range("I136").Value = 1425.2142
$ws.Range("J136").Value = 2719
$ws.Range("K136").Value = 4275.642599999999
$ws.Range("L136").Value = 8157
$ws.Range("M136").Value = -1725.642599999999
$ws.Range("N136").Value = -13257

$ws.Range("H139").Value = 50392
$ws.Range("I139").Value = 21998.4
$ws.Range("J139").Value = 97714.664
$ws.Range("K139").Value = 21998.4
$ws.Range("M139").Value = -16858.4
$ws.Range("N139").Value = -107994.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4371.143
$ws.Range("I62").Value = 4325
$ws.Range("J62").Value = 4432.6665
$ws.Range("K62").Value = 4325
$ws.Range("L62").Value = 4432.6665
$ws.Range("M62").Value = -3701
$ws.Range("N62").Value = -5680.6665

$ws.Range("H65").Value = 4371.143
$ws.Range("I65").Value = 4325
$ws.Range("J65").Value = 4432.6665
$ws.Range("K65").Value = 21625
$ws.Range("L65").Value = 22163.3325
$ws.Range("M65").Value = -18505
$ws.Range("N65").Value = -28403.3325

$ws.Range("H81").Value = 8071.143
$ws.Range("I81").Value = 1090.5454
$ws.Range("J81").Value = 33666.668
$ws.Range("K81").Value = 2181.0908
$ws.Range("L81").Value = 67333.336
$ws.Range("M81").Value = -1120.0908
$ws.Range("N81").Value = -69455.336

$ws.Range("H84").Value = 8071.143
$ws.Range("I84").Value = 1090.5454
$ws.Range("J84").Value = 33666.668
$ws.Range("K84").Value = 10905.454
$ws.Range("L84").Value = 336666.68
$ws.Range("M84").Value = -5601.454
$ws.Range("N84").Value = -347274.68

$ws.Range("H107").Value = 1764.625
$ws.Range("I107").Value = 1191.4706
$ws.Range("K107").Value = 3574.4118
$ws.Range("M107").Value = -1654.4118

$ws.Range("H132").Value = 1698.881
$ws.Range("I132").Value = 1535.6316
$ws.Range("K132").Value = 4606.8948
$ws.Range("M132").Value = -2076.8948

$ws.Range("H139").Value = 79965
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").ClearContents()

